$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(17.9291237013249,10.65515354249863,15.61699124479144,36.1467693177711,3.641039648125945,22.94196388563391,7.569573444631312,13.09413966396252,22.68805403911152)
    3 = @(17.27938774151178,10.42123317474463,15.64609786729766,36.23250906424054,3.643132908467952,23.10283232750196,7.570510956830328,13.04478854700499,22.81214387602837)
    4 = @(16.86827079641165,10.27424903423727,15.66601940597313,36.29526626944796,3.644485713945827,23.2072778779854,7.571423973029622,13.01592461059121,22.89477091529578)
    5 = @(16.69790280325673,10.21356765378022,15.67465341429354,36.32337519979434,3.645054029545829,23.25126705433989,7.571881252325176,13.00453120487604,22.93005487862413)
    6 = @(16.66944924377189,10.20344600352394,15.67611824843167,36.32819551446477,3.645149428465611,23.25865761550288,7.571962341220847,13.0026618251816,22.93601100239907)
    7 = @(16.86598431293348,10.27343375765641,15.66613375826072,36.29563510402616,3.644493309396547,23.2078653536728,7.571429794519095,13.01576945213145,22.8952402443284)
    8 = @(17.70775105687143,10.57522923098505,15.62660195402019,36.174228900454,3.641747420012167,22.9962545016546,7.569826904257767,13.07682946833388,22.7295012714404)
    9 = @(19.253032483317,11.13795973317194,15.56533006809712,36.01671638708815,3.636896126459893,22.62626243978206,7.569344167591739,13.20764129032455,22.45579352087245)
    10 = @(20.3138492535002,11.5306099632594,15.53019803199958,35.95049870907373,3.633653578576229,22.38181245892528,7.570589784105695,13.3100312917006,22.28631208298763)
    11 = @(20.77851123028084,11.7041461225618,15.51635655054137,35.93118910620026,3.632247574311912,22.27654653485375,7.571499277838702,13.35786951540654,22.21614684116147)
    12 = @(20.95176945132978,11.7690853927443,15.51142240063053,35.92543557474305,3.631725030005859,22.23753865309666,7.571892572309132,13.37615645759096,22.19058005978973)
    13 = @(20.91457704065622,11.75513470005,15.5124713975991,35.92660533453916,3.631837130691768,22.24590170520774,7.571805701595537,13.37221055896821,22.19604160158109)
    14 = @(20.79282003590228,11.70950449091812,15.51594445864906,35.93068450585542,3.632204386572437,22.27332021333639,7.571530656990951,13.35937060901212,22.21402330318703)
    15 = @(20.71788544301934,11.68145247399249,15.51811181678166,35.93338618219568,3.63243062646595,22.29022607975326,7.571368538653636,13.35152783287814,22.22516845899834)
    16 = @(20.2831111578823,11.51916283298844,15.53114563551543,35.9519785640221,3.633746849391537,22.38881131429243,7.570537210476787,13.30692949942557,22.29103759394415)
    17 = @(20.01171085044856,11.41826906955696,15.53968932732566,35.96615684968599,3.634571958514206,22.45081072960247,7.570114738819407,13.27988594744259,22.33322633254587)
    18 = @(19.8539291996639,11.35976168912563,15.5448049135279,35.97532945614642,3.635053041449921,22.48702980120863,7.569904040600256,13.26445044696105,22.35814433860338)
    19 = @(19.80022265802365,11.33987175524182,15.54657158056453,35.9786097928337,3.635217046115526,22.4993888916247,7.569838261244362,13.25924501946618,22.36669296265045)
    20 = @(20.04077663499483,11.42905897426833,15.53875898736981,35.96454219587246,3.634483451670419,22.44415296478121,7.570156372075829,13.28275250819203,22.32866772715618)
    21 = @(20.82865718143167,11.72292856118796,15.51491599966842,35.92944403050435,3.632096246900361,22.26524354976376,7.571610120746668,13.36313743292652,22.20871436235601)
    22 = @(21.32780344964625,11.91045330284324,15.50112427963678,35.91559098601471,3.630593630035658,22.15329423198942,7.572844996245652,13.41666948964527,22.13616869212285)
    23 = @(21.06288028597202,11.81079659673272,15.50832145785632,35.92215234309478,3.631390354628738,22.21258796912053,7.572160001689249,13.38801053671219,22.17435023511753)
    24 = @(20.02764143508803,11.42418241927021,15.53917895951833,35.96526899978658,3.634523444688342,22.44716114983148,7.570137449382369,13.28145618587351,22.33072660745455)
    25 = @(18.84741735236711,10.98916925574701,15.58016819594602,36.05065848301183,3.638151787156346,22.72154421189933,7.569192149936386,13.17111481076837,22.52431163288054)
}

$cols = @("B","C","E","F","G","I","J","L","O")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = "$($cols[$i])$r"
        $ws.Range($cell).Value = $vals[$i]
    }
}
